$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# All-night Crafting | Cunning Craftsman's Tisane (row 138)
$ws.Range("H138").Value = 1377.86
$ws.Range("I138").Value = 703.4545000000001
$ws.Range("J138").Value = 1710.0299
$ws.Range("K138").Value = 2110.3635
$ws.Range("L138").Value = 5130.0897
$ws.Range("M138").Value = 3029.6365
$ws.Range("N138").Value = -15410.0897

# Remedy for Reason | Grade 1 Gemdraught of Mind (row 141)
$ws.Range("H141").Value = 872.7778
$ws.Range("I141").Value = 872.7778
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2618.3334
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2561.6666
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Signed, Shield, Delivered | Titanbronze Tower Shield (row 117)
$ws.Range("H117").Value = 51224
$ws.Range("J117").Value = 51224
$ws.Range("L117").Value = 51224
$ws.Range("N117").Value = -60402

# Trial and Error | Dwarven Mythril Chainmail of Fending (row 119)
$ws.Range("H119").Value = 24499
$ws.Range("J119").Value = 24499
$ws.Range("L119").Value = 24499
$ws.Range("N119").Value = -34175

# Haste for High Durium | High Durium Nugget (row 122)
$ws.Range("H122").Value = 1821.0769
$ws.Range("I122").Value = 1971.5454
$ws.Range("J122").Value = 993.5
$ws.Range("K122").Value = 5914.6362
$ws.Range("L122").Value = 2980.5
$ws.Range("M122").Value = -3464.6362
$ws.Range("N122").Value = -7880.5

# Don't Bore Me, Ore Me | Mountain Chromite Ingot (row 132)
$ws.Range("H132").Value = 2245.6191
$ws.Range("I132").Value = 1897.7778
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 5693.3334
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -3163.3334
$ws.Range("N132").Value = -18057.9995

# Shielding My Students | Mountain Chromite Tower Shield (row 133)
$ws.Range("H133").Value = 28121.576
$ws.Range("J133").Value = 28126.44
$ws.Range("L133").Value = 28126.44
$ws.Range("N133").Value = -33186.44

# Forgiveness for My Shins | Ruthenium Sabatons of Fending (row 135)
$ws.Range("H135").Value = 17459.25
$ws.Range("J135").Value = 17459.25
$ws.Range("L135").Value = 17459.25
$ws.Range("N135").Value = -27599.25

$ws = $wb.Worksheets.Item("BSM")
# Ruthenium Supremium | Ruthenium Ingot (row 134)
$ws.Range("H134").Value = 12739.444
$ws.Range("I134").Value = 1832
$ws.Range("J134").Value = 99999
$ws.Range("K134").Value = 5496
$ws.Range("L134").Value = 299997
$ws.Range("M134").Value = -2961
$ws.Range("N134").Value = -305067

$ws = $wb.Worksheets.Item("CRP")
# Wall Not Found | Walnut Lumber (row 31)
$ws.Range("H31").Value = 2210.1035
$ws.Range("I31").Value = 1091.5
$ws.Range("K31").Value = 1091.5
$ws.Range("M31").Value = -796.5

# Armoires of the Rich and Famous | Walnut Lumber (row 34)
$ws.Range("H34").Value = 2210.1035
$ws.Range("I34").Value = 1091.5
$ws.Range("K34").Value = 1091.5
$ws.Range("M34").Value = -889.5

# O Pine | Pine Lumber (row 99)
$ws.Range("H99").Value = 1911.3334
$ws.Range("I99").Value = 1911.3334
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1911.3334
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -413.3334
$ws.Range("N99").ClearContents()

# A Better Conductor | Red Pine Lumber (row 126)
$ws.Range("H126").Value = 1911.3334
$ws.Range("I126").Value = 1911.3334
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5734.0002
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3264.0002
$ws.Range("N126").ClearContents()

# Wood You Be Quiet | Ceiba Lumber (row 134)
$ws.Range("H134").Value = 12346634
$ws.Range("I134").Value = 14493690
$ws.Range("J134").Value = 1057
$ws.Range("K134").Value = 43481070
$ws.Range("L134").Value = 3171
$ws.Range("M134").Value = -43478535
$ws.Range("N134").Value = -8241

$ws = $wb.Worksheets.Item("CUL")
# Such a Butter Face | Fermented Butter (row 68)
$ws.Range("H68").Value = 1176.909
$ws.Range("I68").Value = 798.5
$ws.Range("J68").Value = 1261
$ws.Range("K68").Value = 2395.5
$ws.Range("L68").Value = 3783
$ws.Range("M68").Value = -1584.5
$ws.Range("N68").Value = -5405

# No Margarine of Error (L) | Fermented Butter (row 71)
$ws.Range("H71").Value = 1176.909
$ws.Range("I71").Value = 798.5
$ws.Range("J71").Value = 1261
$ws.Range("K71").Value = 7186.5
$ws.Range("L71").Value = 11349
$ws.Range("M71").Value = -3130.5
$ws.Range("N71").Value = -19461

# Salt of the North | Northern Sea Salt (row 122)
$ws.Range("H122").Value = 443
$ws.Range("I122").Value = 449.33334
$ws.Range("K122").Value = 4044.00006
$ws.Range("M122").Value = -1594.00006

# The Mountain Steeped | Tsai tou Vounou (row 131)
$ws.Range("H131").Value = 25642492
$ws.Range("I131").Value = 200000260
$ws.Range("J131").Value = 1642.6471
$ws.Range("K131").Value = 600000780
$ws.Range("L131").Value = 4927.9413
$ws.Range("M131").Value = -599995740
$ws.Range("N131").Value = -15007.9413

# More Mezcal | Cooking Mezcal (row 132)
$ws.Range("H132").Value = 1995
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1995
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17955
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -23015

# Don't Knock It Till You've Tried It | Mezcal-marinated Swampmonk (row 134)
$ws.Range("H134").Value = 3880.2778
$ws.Range("I134").Value = 938
$ws.Range("J134").Value = 5752.636
$ws.Range("K134").Value = 2814
$ws.Range("L134").Value = 17257.908
$ws.Range("M134").Value = 2256
$ws.Range("N134").Value = -27397.908

$ws = $wb.Worksheets.Item("GSM")
# Needs More Prayerbell | Hardsilver Ingot (row 80)
$ws.Range("H80").Value = 5281.2
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5281.2
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 5281.2
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -7277.2

# With a Noise That Reaches Heaven (L) | Hardsilver Ingot (row 83)
$ws.Range("H83").Value = 5281.2
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5281.2
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 26406
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -36390

# Awarding Academic Excellence | Ametrine (row 122)
$ws.Range("H122").Value = 502160
$ws.Range("I122").Value = 1980
$ws.Range("K122").Value = 5940
$ws.Range("M122").Value = -3490

# Gold Rush Order | Phrygian Gold Ingot (row 126)
$ws.Range("H126").Value = 2444.875
$ws.Range("I126").Value = 1948.25
$ws.Range("J126").Value = 2941.5
$ws.Range("K126").Value = 5844.75
$ws.Range("L126").Value = 8824.5
$ws.Range("M126").Value = -3374.75
$ws.Range("N126").Value = -13764.5

# On Board for Lar | Lar Ingot (row 132)
$ws.Range("H132").Value = 1975.3684
$ws.Range("I132").Value = 1533.375
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 4600.125
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -2070.125
$ws.Range("N132").Value = -18057.9995

$ws = $wb.Worksheets.Item("LTW")
# Tan Before the Ban | Leather (row 7)
$ws.Range("H7").Value = 1856.125
$ws.Range("I7").Value = 1737
$ws.Range("J7").Value = 2213.5
$ws.Range("K7").Value = 1737
$ws.Range("L7").Value = 2213.5
$ws.Range("M7").Value = -1625
$ws.Range("N7").Value = -2437.5

# Saddle Sore | Hard Leather (row 16)
$ws.Range("H16").Value = 1694.9231
$ws.Range("I16").Value = 1728.25
$ws.Range("J16").Value = 1641.6
$ws.Range("K16").Value = 1728.25
$ws.Range("L16").Value = 1641.6
$ws.Range("M16").Value = -1558.25
$ws.Range("N16").Value = -1981.6

# Skin off Their Backs | Aldgoat Leather (row 22)
$ws.Range("H22").Value = 1417.037
$ws.Range("I22").Value = 1605.7142
$ws.Range("J22").Value = 1213.8462
$ws.Range("K22").Value = 1605.7142
$ws.Range("L22").Value = 1213.8462
$ws.Range("M22").Value = -1310.7142
$ws.Range("N22").Value = -1803.8462

# Fire and Hide | Aldgoat Leather (row 27)
$ws.Range("H27").Value = 1417.037
$ws.Range("I27").Value = 1605.7142
$ws.Range("J27").Value = 1213.8462
$ws.Range("K27").Value = 1605.7142
$ws.Range("L27").Value = 1213.8462
$ws.Range("M27").Value = -1498.7142
$ws.Range("N27").Value = -1427.8462

# Best Served Toad | Toad Leather (row 40)
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3272

# Supply Side Logic | Boar Leather (row 46)
$ws.Range("H46").Value = 5571.364
$ws.Range("I46").Value = 2746.25
$ws.Range("J46").Value = 7185.7144
$ws.Range("K46").Value = 2746.25
$ws.Range("L46").Value = 7185.7144
$ws.Range("M46").Value = -2558.25
$ws.Range("N46").Value = -7561.7144

# Battered Books | Saiga Leather (row 126)
$ws.Range("H126").Value = 1856.125
$ws.Range("I126").Value = 1737
$ws.Range("J126").Value = 2213.5
$ws.Range("K126").Value = 5211
$ws.Range("L126").Value = 6640.5
$ws.Range("M126").Value = -2741
$ws.Range("N126").Value = -11580.5

# The Perfect Accessory | Loboskin Amulet of Fending (row 133)
$ws.Range("H133").Value = 39250
$ws.Range("J133").Value = 39250
$ws.Range("L133").Value = 39250
$ws.Range("N133").Value = -44310

# Respect for Br'aax | Br'aax Leather (row 136)
$ws.Range("H136").Value = 1764.7778
$ws.Range("I136").Value = 1438.5
$ws.Range("J136").Value = 2417.3333
$ws.Range("K136").Value = 4315.5
$ws.Range("L136").Value = 7251.999899999999
$ws.Range("M136").Value = -1765.5
$ws.Range("N136").Value = -12351.9999

$ws = $wb.Worksheets.Item("WVR")
# Where the Dragonflies, the Net Catches | Crawler Silk (row 81)
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()

# To Kill a Dragon on Nameday (L) | Crawler Silk (row 84)
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

# Of Great Import | Kudzu Thread (row 100)
$ws.Range("H100").Value = 1860.3334
$ws.Range("I100").Value = 1106.2858
$ws.Range("K100").Value = 2212.5716
$ws.Range("M100").Value = -1671.5716

# Comfy Cabins | Snow Cotton Cloth (row 132)
$ws.Range("H132").Value = 3069.6191
$ws.Range("I132").Value = 2866.5789
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 8599.736699999999
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -6069.736699999999
$ws.Range("N132").Value = -20055.5

# Begin with the Basics | Snow Cotton Jacket (row 133)
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
